# Applies the "Elimina EC anteriores y se agregan nuevos, se modifica base de datos"
# change: replaces the single-worker account-statement table with a
# five-worker table (35 data rows instead of 7) and updates the totals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Make room in the sheet: the table used to have 1 worker x 7
#    periods (rows 16-29, with row 29 using the special "closing"
#    border style). We need 5 workers x 7 periods (rows 16-50, with
#    row 50 using that special closing style). Insert 21 rows right
#    before the old closing row (29) so it is pushed down to become
#    the new closing row (50), and the footer rows below it shift
#    down by the same amount (34->55, 35->56).
# ------------------------------------------------------------------
$ws.Range("B29:J49").Insert()

# ------------------------------------------------------------------
# 2. The freshly inserted rows don't carry the normal data-row
#    border/alignment style (s="15".."20"). Stamp it in by copying
#    the still-correctly-styled first worker's block (B16:J22) over
#    each of the newly-opened 7-row bands.
# ------------------------------------------------------------------
$srcBlock = $ws.Range("B16:J22")
$srcBlock.Copy($ws.Range("B23:J29"))
$srcBlock.Copy($ws.Range("B30:J36"))
$srcBlock.Copy($ws.Range("B37:J43"))
$srcBlock.Copy($ws.Range("B44:J49"))

# ------------------------------------------------------------------
# 3. Write the new worker / period data over the whole 35-row table.
# ------------------------------------------------------------------
$workers = @(
  @{id="9144308"; name="ARNOL DE JESUS MEDINA JIMENEZ"},
  @{id="72217178"; name="EFRAIN EDURADO CHACON GARCIA"},
  @{id="1052943894"; name="LUZ ESTHER BELEÑO SAENZ"},
  @{id="3875968"; name="ELVIS TORRES RENTERIA"},
  @{id="19873281"; name="NICOLAS ROJAS MUÑOZ"}
)
$periods = @("2502","2501","2412","2411","2410","2409","2408")

$row = 16
for ($w = 0; $w -lt $workers.Length; $w++) {
  for ($p = 0; $p -lt $periods.Length; $p++) {
    $ws.Cells.Item($row, 2).Value = "CC"
    $ws.Cells.Item($row, 3).Value = $workers[$w].id
    $ws.Cells.Item($row, 4).Value = $workers[$w].name
    $ws.Cells.Item($row, 5).Value = $periods[$p]
    if ($p -eq 0) {
      $ws.Cells.Item($row, 6).Value = 32933
    } else {
      $ws.Cells.Item($row, 6).Value = 52000
    }
    $ws.Cells.Item($row, 7).Value = 1300000
    $row = $row + 1
  }
}

# ------------------------------------------------------------------
# 4. Update the summary cells: total overdue amount and worker count.
# ------------------------------------------------------------------
$ws.Range("E11").Value = 1724665
$ws.Range("C13").Value = 5
